$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "st_map" (sheet1): add several new fields (columns) and new sample
# rows, add "IsSkipZeroValue" flag to the generator-options JSON, and remove
# the now-unused "StData/" path column.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("st_map")

# --- clear the handful of old cells that have no counterpart in the new
#     layout (keeps everything else, including the header comments / vml
#     drawing relationship, intact) ---
$ws.Range("D1").ClearContents() | Out-Null
$ws.Range("B4").ClearContents() | Out-Null
$ws.Range("C6").ClearContents() | Out-Null
$ws.Range("D4").Clear() | Out-Null
$ws.Range("D7").Clear() | Out-Null

# --- row 1 : class meta ---
$ws.Range("A1").Value = "st_map"
$ws.Range("B1").Value = "st_mapTable"
$ws.Range("C1").Value = '{"IsSkipZeroValue":false,"IsStringId":false,"IsGenItemClass":true,"JSONName":"st_mapJSON","IsGenEnum":false}'

# --- row 2 : field names/types ---
$ws.Range("A2").Value = "id"
$ws.Range("B2").Value = "row"
$ws.Range("C2").Value = "col"
$ws.Range("D2").Value = "test:float"
$ws.Range("E2").Value = "map"
$ws.Range("F2").Value = "tesbo:bool"
$ws.Range("G2").Value = "a0"
$ws.Range("H2").Value = "a1"
$ws.Range("I2").Value = "a2"
$ws.Range("J2").Value = "b0:float"
$ws.Range("K2").Value = "b1"
$ws.Range("L2").Value = "b2"
$ws.Range("M2").Value = "c0:bool"
$ws.Range("N2").Value = "c1"
$ws.Range("O2").Value = "c2"
$ws.Range("P2").Value = "s0:string"
$ws.Range("Q2").Value = "s1"
$ws.Range("R2").Value = "s2"

# --- row 3 ---
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = "'a10101010101010100"
$ws.Range("F3").Value = "ssd"
$ws.Range("G3").Value = 6
$ws.Range("J3").Value = 4
$ws.Range("M3").Value = "ssd"
$ws.Range("N3").Value = "f"
$ws.Range("P3").Value = "dfdf"
$ws.Range("Q3").Value = "dfdf"

# --- row 4 ---
$ws.Range("A4").Value = 1
$ws.Range("C4").Value = 0
# E4 keeps the quote-prefix text style but carries no value: copy the
# format from a cell that already has it, then drop the contents.
$ws.Range("E3").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null
$ws.Range("E4").ClearContents() | Out-Null
$ws.Range("F4").Value = $true
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 5
$ws.Range("I4").Value = 6
$ws.Range("J4").Value = 7.6
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = $true
$ws.Range("N4").Value = $false
$ws.Range("P4").Value = "dfdf"

# --- row 5 ---
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = "'0.0"
$ws.Range("E3").Copy() | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null
$ws.Range("E5").ClearContents() | Out-Null
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 7

# --- row 6 ---
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 0
$ws.Range("D6").Value = 1.5
$ws.Range("E3").Copy() | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").ClearContents() | Out-Null
$ws.Range("F6").Value = $false
$ws.Range("G6").Value = 4
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = 4
$ws.Range("M6").Value = $false
$ws.Range("N6").Value = $false
$ws.Range("O6").Value = $false
$ws.Range("P6").Value = "dfdf"
$ws.Range("Q6").Value = "dfdf"
$ws.Range("R6").Value = "df"

# --- row 7 ---
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 4
$ws.Range("E7").Value = "'a10101010101010100"
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = 6
$ws.Range("L7").Value = 7
$ws.Range("M7").Value = 0
$ws.Range("P7").Value = "dfdf"
$ws.Range("Q7").Value = "dfdf"
$ws.Range("R7").Value = "df"

# --- row 8 (new row) ---
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = 4
$ws.Range("C8").Value = 4
$ws.Range("E8").Value = "'a10101010101010100"
$ws.Range("F8").Value = "'false"
$ws.Range("G8").Value = 5
$ws.Range("H8").Value = 6
$ws.Range("M8").Value = "'false"
$ws.Range("P8").Value = "dfdf"
$ws.Range("Q8").Value = "dfdf"

# ---------------------------------------------------------------------------
# View state: st_map becomes the active sheet/tab, Sheet2 loses focus.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Sheet2")
$ws4.Activate()
$ws4.Range("B11").Select() | Out-Null

$ws.Activate()
$ws.Range("I6").Select() | Out-Null
